$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 22

# Copy formatting (style) from the previous last data row so the new
# row's cells share the same style index as the rest of the table.
$ws.Range("A21:T21").Copy()
$ws.Range("A22:T22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item($newRow, 1).Value = 84645000
$ws.Cells.Item($newRow, 2).Value = "Чойский муниципальный район "
$ws.Cells.Item($newRow, 3).Value = 2018
$ws.Cells.Item($newRow, 4).Value = -0.13492063492063491
$ws.Cells.Item($newRow, 5).Value = [double]"5.5531438875346437E-2"
$ws.Cells.Item($newRow, 6).Value = 0.34933223314617218
$ws.Cells.Item($newRow, 7).Value = 0.42370510470138839
$ws.Cells.Item($newRow, 8).Value = 0.37399784084144988
$ws.Cells.Item($newRow, 9).Value = 0.20103175623387251
$ws.Cells.Item($newRow, 10).Value = 0.41700794427406851
$ws.Cells.Item($newRow, 11).Value = 0.32049206863062479
$ws.Cells.Item($newRow, 12).Value = 0.19370555070779899
$ws.Cells.Item($newRow, 13).Value = 0.1046713062855934
$ws.Cells.Item($newRow, 14).Value = 0.23519158446383381
$ws.Cells.Item($newRow, 15).Value = [double]"6.5909613799584377E-2"
$ws.Cells.Item($newRow, 16).Value = 0.11436281212698141
$ws.Cells.Item($newRow, 17).Value = [double]"5.5137836962902788E-2"
$ws.Cells.Item($newRow, 18).Value = 0.31581479254359662
$ws.Cells.Item($newRow, 19).Value = 0.58127579720865963
$ws.Cells.Item($newRow, 20).Value = [double]"3.6359010741574729E-2"
